# Apply the updated crypto price/volume figures scraped for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.415.76"
$ws.Range("E2").Value = "  +12.75%  "
$ws.Range("D3").Value = "1.819.11"
$ws.Range("E3").Value = "  +7.78%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.66%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.62"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("E10").Value = "  +6.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0683"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0931"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.36%  "
$ws.Range("D13").Value = "2.079.69"
$ws.Range("E13").Value = "  +7.68%  "
$ws.Range("D14").Value = "1.820.19"
$ws.Range("E14").Value = "  +6.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.647"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.38%  "
$ws.Range("D16").Value = "34.378.80"
$ws.Range("E16").Value = "  +12.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "261.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.20%  "
$ws.Range("D21").Value = "0.0₃0754"
$ws.Range("E21").Value = "  +4.75%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.52%  "
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.66%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0517"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.27%  "
$ws.Range("E33").Value = "  +6.45%  "
$ws.Range("E34").Value = "  +8.69%  "
$ws.Range("D35").Value = "1.588.83"
$ws.Range("E35").Value = "  +5.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.58%  "
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "86.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.34%  "
$ws.Range("E39").Value = "  +7.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0189"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.920"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.85%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0523"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.54%  "
$ws.Range("D47").Value = "1.979.68"
$ws.Range("E47").Value = "  +8.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +18.54%  "
